# "add and search counterparty"
# The AddOpportunity sheet has a "Search" helper column (AG) used to look up
# a counterparty. Update the searched-for counterparty from "Tec Alliance"
# to "TEC Canada" on both sample rows, then tidy up the leftover blank
# template rows below the data (rows 5-9 only carried stray number-format
# styling with no values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddOpportunity")

$ws.Range("AG2").Value = "TEC Canada"
$ws.Range("AG3").Value = "TEC Canada"

# Remove the now-unused blank formatted rows trailing the sample data.
$ws.Range("A5:AI9").Clear()

# Leave the cursor where the user last left it.
[void]$ws.Range("AG4").Select()
